# Generate Report for Handback
# Refresh the timestamps recorded for the c9bda275-... handback report:
#   - Overview sheet: "Latest HO Xliff Generate Date" for the file
#   - zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
#   - de-de sheet: "Correspond Handback DateTime"

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G3").Value = "2016-08-27 20:44:59"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H3").Value = "2016-08-27 20:44:55"
$zhcn.Range("K3").Value = "2016-08-27 20:45:30"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("K3").Value = "2016-08-27 20:45:37"
